function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws.Range("D2") "245.12"
Set-TextValue $ws.Range("E2") "-0.59%"

# Row 3
Set-TextValue $ws.Range("D3") "29.13"
Set-TextValue $ws.Range("E3") "-1.64%"

# Row 4
Set-TextValue $ws.Range("D4") "5.260"
Set-TextValue $ws.Range("E4") "1.96%"

# Row 5
Set-TextValue $ws.Range("D5") "0.05705"
Set-TextValue $ws.Range("E5") "-0.03%"

# Row 6
Set-TextValue $ws.Range("D6") "6.614"

# Row 7
Set-TextValue $ws.Range("D7") "3.191"
Set-TextValue $ws.Range("E7") "3.84%"

# Row 8
Set-TextValue $ws.Range("D8") "0.8502"
Set-TextValue $ws.Range("E8") "-0.66%"

# Row 9
Set-TextValue $ws.Range("D9") "0.8585"
Set-TextValue $ws.Range("E9") "-1.35%"

# Row 10
Set-TextValue $ws.Range("D10") "0.1371"
Set-TextValue $ws.Range("E10") "0.33%"

# Row 11
Set-TextValue $ws.Range("D11") "0.07036"
Set-TextValue $ws.Range("E11") "-0.58%"

# Row 12
Set-TextValue $ws.Range("D12") "0.03189"
Set-TextValue $ws.Range("E12") "9.17%"

# Row 13
Set-TextValue $ws.Range("D13") "0.09281"
Set-TextValue $ws.Range("E13") "-1.18%"

# Row 14
Set-TextValue $ws.Range("E14") "0.68%"

# Row 15
Set-TextValue $ws.Range("D15") "0.0005979"
Set-TextValue $ws.Range("E15") "-0.36%"

# Row 16
Set-TextValue $ws.Range("D16") "0.005979"
Set-TextValue $ws.Range("E16") "-1.12%"

# Row 17
Set-TextValue $ws.Range("D17") "3.494"
Set-TextValue $ws.Range("E17") "0.23%"

# Row 18
Set-TextValue $ws.Range("E18") "-4.26%"

# Row 19
Set-TextValue $ws.Range("D19") "0.3160"
Set-TextValue $ws.Range("E19") "-0.46%"

# Row 20
Set-TextValue $ws.Range("D20") "0.03339"
Set-TextValue $ws.Range("E20") "1.02%"

# Row 21
Set-TextValue $ws.Range("E21") "-1.97%"

# Row 22
Set-TextValue $ws.Range("D22") "3.488"
Set-TextValue $ws.Range("E22") "0.67%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04085"
Set-TextValue $ws.Range("E23") "-2.21%"

# Row 24
Set-TextValue $ws.Range("E24") "-0.04%"

# Row 25
Set-TextValue $ws.Range("D25") "0.001221"
Set-TextValue $ws.Range("E25") "-0.07%"

# Row 26
Set-TextValue $ws.Range("D26") "0.004144"
Set-TextValue $ws.Range("E26") "-17.57%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0001200"
Set-TextValue $ws.Range("E27") "-0.78%"

# Row 28
Set-TextValue $ws.Range("D28") "0.0001449"
Set-TextValue $ws.Range("E28") "-25.27%"

# Row 40
Set-TextValue $ws.Range("D40") "0.03758"
Set-TextValue $ws.Range("E40") "0.37%"

# Row 41
Set-TextValue $ws.Range("D41") "0.1064"
Set-TextValue $ws.Range("E41") "-0.67%"

# Row 42
Set-TextValue $ws.Range("D42") "0.003719"
Set-TextValue $ws.Range("E42") "-35.43%"

# Row 43
Set-TextValue $ws.Range("D43") "0.002401"
Set-TextValue $ws.Range("E43") "20.09%"

# Row 44
Set-TextValue $ws.Range("D44") "0.009354"
Set-TextValue $ws.Range("E44") "-4.78%"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005261"
Set-TextValue $ws.Range("E45") "0.94%"

# Row 46
Set-TextValue $ws.Range("E46") "0.00%"

# Row 47
Set-TextValue $ws.Range("D47") "0.07499"
Set-TextValue $ws.Range("E47") "29.27%"

# Row 48
Set-TextValue $ws.Range("E48") "-5.11%"

# Row 49
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "0.00%"

# Row 50
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "0.00%"
